# Applies the "updated figures and maturity info" change to the
# Value_data worksheet:
#  - Inserts a new column C "Market_value_KES" (numeric KES values)
#    shifting the old "Families" / "Common_name" columns right by one.
#  - Renames B1 header from "Market_value" to "Market_value_USD".
#  - Appends three new rows (Gerreidae / Clupeidae / Trichiuridae) with
#    maturity/market-value figures but no Category / Common_name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Value_data")

# Insert a new column before the current column C ("Families"), which
# shifts Families -> D and Common_name -> E.
$ws.Columns("C").Insert()

# Both value columns (USD, KES) share the same width (raw OOXML width of
# 18 "characters" -- ColumnWidth via COM reports ~5/6 wider than the
# stored column width attribute, so compensate for that offset here).
$ws.Range("B:C").ColumnWidth = 18 - (5/6)

# Append three new rows for additional families (no Category / Common_name)
# first, so new shared strings are introduced in the same order as the
# source edit (families before the renamed headers).
$ws.Cells.Item(22, 2).Value = 2
$ws.Cells.Item(22, 3).Value = 250
$ws.Cells.Item(22, 4).Value = "Gerreidae"

$ws.Cells.Item(23, 2).Value = 2
$ws.Cells.Item(23, 3).Value = 250
$ws.Cells.Item(23, 4).Value = "Clupeidae"

$ws.Cells.Item(24, 2).Value = 2
$ws.Cells.Item(24, 3).Value = 250
$ws.Cells.Item(24, 4).Value = "Trichiuridae"

# Update headers on row 1.
$ws.Range("B1").Value = "Market_value_USD"
$ws.Range("C1").Value = "Market_value_KES"

# Fill in the new KES values (column C) for the existing data rows.
$kesValues = @{
    2  = 350
    3  = 350
    4  = 350
    5  = 350
    6  = 250
    7  = 250
    8  = 300
    9  = 300
    10 = 300
    11 = 300
    12 = 300
    13 = 300
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
}

foreach ($row in $kesValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $kesValues[$row]
}

# Match the final selection noted in the saved file.
$ws.Range("C20").Select()
